$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 35: "100k" resistor list cell (C35) ---
# Original: R11,R14,R17,R20,R35,R36,R37,R38,R48,R49,R55,R56
# New:      R11,R14,R17,R20,R35,R36,R37,R38,R48,R49,R55
# (R56 removed - component no longer needed)
$ws.Range("C35").Value = "R11,R14,R17,R20,R35,R36,R37,R38,R48,R49,R55"

# --- Row 30: "1k" resistor list cell (C30) ---
# Original: R10,R13,R16,R19,R21,R23,R24,R29,R30,R39,R50,R51,R57,R58,R59,R62,R64
# New:      R10,R13,R16,R19,R21,R23,R24,R29,R30,R39,R50,R51,R57,R59,R62,R64
# (R58 removed - component no longer needed)
$c30 = $ws.Range("C30")
$c30.Value = "R10,R13,R16,R19,R21,R23,R24,R29,R30,R39,R50,R51,R57,R59,R62,R64"

$c30text = $c30.Value()

# Re-apply the original rich-text coloring that existed on this cell:
#  - "R39" in green (00B050)
#  - ",R50,R51,R57," in black (000000)
#  - "R59" in green (00B050)
#  - ",R62," in black (000000)
#  - "R64" in red (FF0000)
$i = $c30text.IndexOf("R39")
$c30.Characters($i + 1, 3).Font.Color = 5287936

$i = $c30text.IndexOf(",R50,R51,R57,")
$c30.Characters($i + 1, 13).Font.Color = 0

$i = $c30text.IndexOf("R59")
$c30.Characters($i + 1, 3).Font.Color = 5287936

$i = $c30text.IndexOf(",R62,")
$c30.Characters($i + 1, 5).Font.Color = 0

$i = $c30text.IndexOf("R64")
$c30.Characters($i + 1, 3).Font.Color = 255

# Update the view selection to match where the edit was made
$ws.Range("C30").Select()
